$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 507 and 508 (shifts existing rows 507-598 down to 509-600)
$ws.Rows.Item(507).Insert()
$ws.Rows.Item(508).Insert()

# Row 507
$ws.Range("A507").Value = 11
$ws.Range("B507").Value = "Vega Monumental Concepción"
$ws.Range("C507").Value = "Bíobío"
$ws.Range("D507").Value = 45154
$ws.Range("E507").Value = 8
$ws.Range("F507").Value = 100112002
$ws.Range("G507").Value = "Pimiento"
$ws.Range("H507").Value = "Zafiro rojo"
$ws.Range("I507").Value = "Primera"
$ws.Range("J507").Value = 100
$ws.Range("K507").Value = 18000
$ws.Range("L507").Value = 18000
$ws.Range("M507").Value = 18000
$ws.Range("N507").Value = "$/caja 15 kilos"
$ws.Range("O507").Value = "Región de Arica y Parinacota"
$ws.Range("P507").Value = 1200
$ws.Range("Q507").Value = 15
$ws.Range("R507").Value = "Hortaliza"

# Row 508
$ws.Range("A508").Value = 11
$ws.Range("B508").Value = "Vega Monumental Concepción"
$ws.Range("C508").Value = "Bíobío"
$ws.Range("D508").Value = 45154
$ws.Range("E508").Value = 8
$ws.Range("F508").Value = 100112002
$ws.Range("G508").Value = "Pimiento"
$ws.Range("H508").Value = "Zafiro verde"
$ws.Range("I508").Value = "Primera"
$ws.Range("J508").Value = 100
$ws.Range("K508").Value = 17000
$ws.Range("L508").Value = 17000
$ws.Range("M508").Value = 17000
$ws.Range("N508").Value = "$/caja 15 kilos"
$ws.Range("O508").Value = "Región de Arica y Parinacota"
$ws.Range("P508").Value = 1133
$ws.Range("Q508").Value = 15
$ws.Range("R508").Value = "Hortaliza"
